$wb = $excel.ActiveWorkbook

$info = $wb.Worksheets.Item("INFO")
$demands = $wb.Worksheets.Item("Demands")
$sim = $wb.Worksheets.Item("Simulation")

# ---------------------------------------------------------------------------
# 1. Restructure the "Simulation" sheet:
#    - insert a new "timestep" row right after "stop_time" (old row 6)
#    - insert three new rows (p_nominal, T_nominal, fraction_glycol) right
#      after "medium" (old row 8, new row 9)
# ---------------------------------------------------------------------------

# Insert one row at row 6 -> old rows 6,7,8,9 shift to 7,8,9,10
$sim.Rows.Item(6).Insert()
$sim.Rows.Item(6).RowHeight = 18.75

# Insert three rows at row 10 -> old row 9 (now 10) shifts further down to 13
$sim.Rows("10:12").Insert()
$sim.Rows("10:12").RowHeight = 18.75

# --- Fill the new row 6: timestep -----------------------------------------
$sim.Range("A6").Value = "timestep"
$sim.Range("B6").Value = 900
$sim.Range("C6").Value = "s"
$sim.Range("D6").Value = ""

# --- Fill the new rows 10-12: p_nominal / T_nominal / fraction_glycol -----
$sim.Range("A10").Value = "p_nominal"
$sim.Range("B10").Value = 400000
$sim.Range("C10").Value = "Pa"
$sim.Range("D10").Value = "Nominal pressure of medium"

$sim.Range("A11").Value = "T_nominal"
$sim.Range("B11").Value = 353
$sim.Range("C11").Value = "K"
$sim.Range("D11").Value = "Nominal temperature of medium"

$sim.Range("A12").Value = "fraction_glycol"
$sim.Range("B12").Value = 0.3
$sim.Range("C12").Value = "-"
$sim.Range("D12").Value = "Fraction of glycol if water-glycol medium is used"

# --- Apply formatting matching the rest of the workbook --------------------
# Style "2" (plain, no border) is used on A/C/D of the new rows; style "20"
# (number format, no border) is used on the numeric B cells. Copy formats
# from existing cells that already carry these exact styles so no new
# style entries are introduced.
$info.Range("A1").Copy()
$sim.Range("A6").PasteSpecial(-4122)
$sim.Range("C6").PasteSpecial(-4122)
$sim.Range("D6").PasteSpecial(-4122)
$sim.Range("A10:A12").PasteSpecial(-4122)
$sim.Range("C10:C12").PasteSpecial(-4122)
$sim.Range("D10:D12").PasteSpecial(-4122)

$demands.Range("B3").Copy()
$sim.Range("B6").PasteSpecial(-4122)
$sim.Range("B10:B12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Make "Simulation" the active sheet (it was "Pipes" before).
# ---------------------------------------------------------------------------
$sim.Activate()
